# Apply updated cryptocurrency price/volume data pulled on 2023-08-20.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.310.18'
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').Value = '1.680.65'
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').Value = '''218.02'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').Value = '''0.5507'
$ws.Range('E6').Value = '  +7.86%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').Value = '''0.2702'
$ws.Range('E8').Value = '  +1.67%  '
$ws.Range('E9').Value = '  +1.02%  '
$ws.Range('D10').Value = '''22.05'
$ws.Range('E10').Value = '  +0.88%  '
$ws.Range('D11').Value = '''0.07540'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.702.38'
$ws.Range('E12').Value = '  +1.93%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''4.536'
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('D14').Value = '''0.5802'
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('D15').Value = '''0.000008438'
$ws.Range('E15').Value = '  -1.73%  '
$ws.Range('D16').Value = '''64.93'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').Value = '26.343.38'
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('D18').Value = '''4.919'
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('D21').Value = '''191.23'
$ws.Range('E21').Value = '  -0.94%  '
$ws.Range('D22').Value = '''6.217'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').Value = '''146.59'
$ws.Range('E24').Value = '  +1.18%  '
$ws.Range('D25').Value = '''0.1312'
$ws.Range('E25').Value = '  +9.76%  '
$ws.Range('D26').Value = '''7.873'
$ws.Range('E26').Value = '  +3.36%  '
$ws.Range('D27').Value = '''15.80'
$ws.Range('E27').Value = '  +0.54%  '
$ws.Range('D28').Value = '''0.06368'
$ws.Range('E28').Value = '  -0.94%  '
$ws.Range('D29').Value = '''1.396'
$ws.Range('E29').Value = '  +4.45%  '
$ws.Range('D31').Value = '''3.585'
$ws.Range('E31').Value = '  +0.90%  '
$ws.Range('D32').Value = '''3.575'
$ws.Range('E32').Value = '  +1.58%  '
$ws.Range('D33').Value = '''1.671'
$ws.Range('E33').Value = '  +1.22%  '
$ws.Range('E34').Value = '  +1.47%  '
$ws.Range('D35').Value = '''0.6168'
$ws.Range('E35').Value = '  +0.97%  '
$ws.Range('E36').Value = '  +1.31%  '
$ws.Range('D37').Value = '''2.721'
$ws.Range('E37').Value = '  +1.20%  '
$ws.Range('D38').Value = '''6.240'
$ws.Range('E38').Value = '  -0.32%  '
$ws.Range('D39').Value = '1.113.17'
$ws.Range('E39').Value = '  +1.58%  '
$ws.Range('D40').Value = '''0.01624'
$ws.Range('E40').Value = '  +1.19%  '
$ws.Range('D41').Value = '''0.8718'
$ws.Range('E41').Value = '  +1.26%  '
$ws.Range('E42').Value = '  +0.63%  '
$ws.Range('D43').Value = '''100.65'
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').Value = '1.830.54'
$ws.Range('E44').Value = '  +0.71%  '
$ws.Range('E45').Value = '  -5.57%  '
$ws.Range('D46').Value = '''57.35'
$ws.Range('E46').Value = '  +1.62%  '
$ws.Range('D47').Value = '''8.190'
$ws.Range('E47').Value = '  +1.30%  '
$ws.Range('D48').Value = '''0.9994'
$ws.Range('E48').Value = '  -1.00%  '
$ws.Range('D49').Value = '''0.05275'
$ws.Range('E49').Value = '  +0.66%  '
$ws.Range('D50').Value = '''0.4290'
$ws.Range('E50').Value = '  +0.12%  '
$ws.Range('D51').Value = '''6.051'
$ws.Range('E51').Value = '  -0.01%  '
